# Updates the loading-percent result values for the "380 kV" case (Case_5_84).
# Each data row (r=2..25) gets new computed values in columns B-F, H-L, N-O
# (columns G and M remain 0 and are left untouched).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 13.99924139758921
$ws.Range("C2").Value = 6.361664125099642
$ws.Range("D2").Value = 8.045919013763553
$ws.Range("E2").Value = 12.61410445800589
$ws.Range("F2").Value = 36.2323774264139
$ws.Range("H2").Value = 7.344005520526261
$ws.Range("I2").Value = 27.12140291937764
$ws.Range("J2").Value = 9.893607410861501
$ws.Range("K2").Value = 11.25973831386665
$ws.Range("L2").Value = 11.15898360159033
$ws.Range("N2").Value = 20.53014082243464
$ws.Range("O2").Value = 27.93921226684159

# Row 3
$ws.Range("B3").Value = 13.80182180899699
$ws.Range("C3").Value = 6.298212151884246
$ws.Range("D3").Value = 8.029415387163716
$ws.Range("E3").Value = 12.63012674708947
$ws.Range("F3").Value = 36.31096809774417
$ws.Range("H3").Value = 7.344005520526261
$ws.Range("I3").Value = 27.2078781511722
$ws.Range("J3").Value = 9.911393326257361
$ws.Range("K3").Value = 11.12130971306556
$ws.Range("L3").Value = 11.15848792935492
$ws.Range("N3").Value = 20.58736145625739
$ws.Range("O3").Value = 28.0181076383436

# Row 4
$ws.Range("B4").Value = 13.68148359534498
$ws.Range("C4").Value = 6.258371957167349
$ws.Range("D4").Value = 8.020352362113185
$ws.Range("E4").Value = 12.64144690360763
$ws.Range("F4").Value = 36.36595609200749
$ws.Range("H4").Value = 7.344005520526261
$ws.Range("I4").Value = 27.26507884983396
$ws.Range("J4").Value = 9.922986257986219
$ws.Range("K4").Value = 11.03714018040991
$ws.Range("L4").Value = 11.15955070620223
$ws.Range("N4").Value = 20.6241636159496
$ws.Range("O4").Value = 28.07121664278721

# Row 5
$ws.Range("B5").Value = 13.63272346479568
$ws.Range("C5").Value = 6.241922560646859
$ws.Range("D5").Value = 8.016930988503656
$ws.Range("E5").Value = 12.64643330601867
$ws.Range("F5").Value = 36.39005591659537
$ws.Range("H5").Value = 7.344005520526261
$ws.Range("I5").Value = 27.28942094007816
$ws.Range("J5").Value = 9.927879950347396
$ws.Range("K5").Value = 11.00308519917003
$ws.Range("L5").Value = 11.16032879750205
$ws.Range("N5").Value = 20.63958154309439
$ws.Range("O5").Value = 28.09403195708266

# Row 6
$ws.Range("B6").Value = 13.62464540045794
$ws.Range("C6").Value = 6.239178426502607
$ws.Range("D6").Value = 8.016379370010773
$ws.Range("E6").Value = 12.64728385996028
$ws.Range("F6").Value = 36.39415981033879
$ws.Range("H6").Value = 7.344005520526261
$ws.Range("I6").Value = 27.29352528300624
$ws.Range("J6").Value = 9.928702792736271
$ws.Range("K6").Value = 10.99744625574508
$ws.Range("L6").Value = 11.16047887015136
$ws.Range("N6").Value = 20.64216712565747
$ws.Range("O6").Value = 28.09789124916151

# Row 7
$ws.Range("B7").Value = 13.68082479454695
$ws.Range("C7").Value = 6.2581509725288
$ws.Range("D7").Value = 8.0203051159181
$ws.Range("E7").Value = 12.64151263955469
$ws.Range("F7").Value = 36.36627426264709
$ws.Range("H7").Value = 7.344005520526261
$ws.Range("I7").Value = 27.26540295557463
$ws.Range("J7").Value = 9.923051569258066
$ws.Range("K7").Value = 11.03667986320947
$ws.Range("L7").Value = 11.15955980136432
$ws.Range("N7").Value = 20.62436984231998
$ws.Range("O7").Value = 28.07151958960301

# Row 8
$ws.Range("B8").Value = 13.9310235597562
$ws.Range("C8").Value = 6.339972895055588
$ws.Range("D8").Value = 8.040008133824871
$ws.Range("E8").Value = 12.61932162937032
$ws.Range("F8").Value = 36.25807734114375
$ws.Range("H8").Value = 7.344005520526261
$ws.Range("I8").Value = 27.15036781110021
$ws.Range("J8").Value = 9.899600710571562
$ws.Range("K8").Value = 11.21185895822718
$ws.Range("L8").Value = 11.15852971036876
$ws.Range("N8").Value = 20.54952493412331
$ws.Range("O8").Value = 27.96544644021783

# Row 9
$ws.Range("B9").Value = 14.42609578587068
$ws.Range("C9").Value = 6.493131670590581
$ws.Range("D9").Value = 8.087009822689211
$ws.Range("E9").Value = 12.58754101516121
$ws.Range("F9").Value = 36.09937763561176
$ws.Range("H9").Value = 7.344005520526261
$ws.Range("I9").Value = 26.9573421680558
$ws.Range("J9").Value = 9.858929069082972
$ws.Range("K9").Value = 11.56030456700002
$ws.Range("L9").Value = 11.16730044604803
$ws.Range("N9").Value = 20.4159370847742
$ws.Range("O9").Value = 27.79449176712506

# Row 10
$ws.Range("B10").Value = 14.78917022560758
$ws.Range("C10").Value = 6.600811010715073
$ws.Range("D10").Value = 8.126459528955239
$ws.Range("E10").Value = 12.57130983079086
$ws.Range("F10").Value = 36.01543850277299
$ws.Range("H10").Value = 7.344005520526261
$ws.Range("I10").Value = 26.83536054131724
$ws.Range("J10").Value = 9.832261652052688
$ws.Range("K10").Value = 11.81712152375048
$ws.Range("L10").Value = 11.18024224746065
$ws.Range("N10").Value = 20.32574744413262
$ws.Range("O10").Value = 27.69151239848753

# Row 11
$ws.Range("B11").Value = 14.95350794583671
$ws.Range("C11").Value = 6.648663129478225
$ws.Range("D11").Value = 8.145431851618554
$ws.Range("E11").Value = 12.56546295852031
$ws.Range("F11").Value = 35.98435089938584
$ws.Range("H11").Value = 7.344005520526261
$ws.Range("I11").Value = 26.7841720196384
$ws.Range("J11").Value = 9.820822272756704
$ws.Range("K11").Value = 11.9336743279518
$ws.Range("L11").Value = 11.18752162232018
$ws.Range("N11").Value = 20.28642881265313
$ws.Range("O11").Value = 27.64958219251858

# Row 12
$ws.Range("B12").Value = 15.01556370951038
$ws.Range("C12").Value = 6.666614043326772
$ws.Range("D12").Value = 8.152759928950349
$ws.Range("E12").Value = 12.56346909467981
$ws.Range("F12").Value = 35.97359939923598
$ws.Range("H12").Value = 7.344005520526261
$ws.Range("I12").Value = 26.76540668667843
$ws.Range("J12").Value = 9.816589528472925
$ws.Range("K12").Value = 11.97773307597781
$ws.Range("L12").Value = 11.19047649042017
$ws.Range("N12").Value = 20.27178439067902
$ws.Range("O12").Value = 27.63441159304872

# Row 13
$ws.Range("B13").Value = 15.00220756894009
$ws.Range("C13").Value = 6.662755645579077
$ws.Range("D13").Value = 8.151175370864157
$ws.Range("E13").Value = 12.56388872728074
$ws.Range("F13").Value = 35.97586953137446
$ws.Range("H13").Value = 7.344005520526261
$ws.Range("I13").Value = 26.76942062077646
$ws.Range("J13").Value = 9.817496724091807
$ws.Range("K13").Value = 11.96824826441023
$ws.Range("L13").Value = 11.18983131670961
$ws.Range("N13").Value = 20.27492746329231
$ws.Range("O13").Value = 27.63764737983387

# Row 14
$ws.Range("B14").Value = 14.95861714220864
$ws.Range("C14").Value = 6.650143400711168
$ws.Range("D14").Value = 8.146031880656727
$ws.Range("E14").Value = 12.56529451315821
$ws.Range("F14").Value = 35.98344591060278
$ws.Range("H14").Value = 7.344005520526261
$ws.Range("I14").Value = 26.78261578520328
$ws.Range("J14").Value = 9.820472058347685
$ws.Range("K14").Value = 11.9373008293455
$ws.Range("L14").Value = 11.18776075732548
$ws.Range("N14").Value = 20.28521910939808
$ws.Range("O14").Value = 27.64831991644391

# Row 15
$ws.Range("B15").Value = 14.93189226543229
$ws.Range("C15").Value = 6.642395742408344
$ws.Range("D15").Value = 8.142899935161326
$ws.Range("E15").Value = 12.5661842532084
$ws.Range("F15").Value = 35.98821958759653
$ws.Range("H15").Value = 7.344005520526261
$ws.Range("I15").Value = 26.79077877364103
$ws.Range("J15").Value = 9.822307431244234
$ws.Range("K15").Value = 11.91833346353717
$ws.Range("L15").Value = 11.18651825110405
$ws.Range("N15").Value = 20.29155487402128
$ws.Range("O15").Value = 27.65494930198456

# Row 16
$ws.Range("B16").Value = 14.77840886285238
$ws.Range("C16").Value = 6.597660394053187
$ws.Range("D16").Value = 8.125239965374353
$ws.Range("E16").Value = 12.57172279697084
$ws.Range("F16").Value = 36.01761297924095
$ws.Range("H16").Value = 7.344005520526261
$ws.Range("I16").Value = 26.83879239742061
$ws.Range("J16").Value = 9.833023130348289
$ws.Range("K16").Value = 11.80949574802269
$ws.Range("L16").Value = 11.17979437951206
$ws.Range("N16").Value = 20.32835130830458
$ws.Range("O16").Value = 27.69435160373715

# Row 17
$ws.Range("B17").Value = 14.68399926080099
$ws.Range("C17").Value = 6.569921865461533
$ws.Range("D17").Value = 8.114666297803602
$ws.Range("E17").Value = 12.57551358774923
$ws.Range("F17").Value = 36.03746275569918
$ws.Range("H17").Value = 7.344005520526261
$ws.Range("I17").Value = 26.86934899233867
$ws.Range("J17").Value = 9.839773774186162
$ws.Range("K17").Value = 11.74262984626281
$ws.Range("L17").Value = 11.17602476366889
$ws.Range("N17").Value = 20.35136169952451
$ws.Range("O17").Value = 27.71978309667517

# Row 18
$ws.Range("B18").Value = 14.62962281444251
$ws.Range("C18").Value = 6.553861272011763
$ws.Range("D18").Value = 8.108681464106819
$ws.Range("E18").Value = 12.5778386475295
$ws.Range("F18").Value = 36.04954780730768
$ws.Range("H18").Value = 7.344005520526261
$ws.Range("I18").Value = 26.88732921228938
$ws.Range("J18").Value = 9.843721700070628
$ws.Range("K18").Value = 11.70414671730709
$ws.Range("L18").Value = 11.17398772759194
$ws.Range("N18").Value = 20.36475759896237
$ws.Range("O18").Value = 27.73487326909857

# Row 19
$ws.Range("B19").Value = 14.61120083330201
$ws.Range("C19").Value = 6.548405417370621
$ws.Range("D19").Value = 8.106671855306908
$ws.Range("E19").Value = 12.57865074732868
$ws.Range("F19").Value = 36.05375431409635
$ws.Range("H19").Value = 7.344005520526261
$ws.Range("I19").Value = 26.89348654062731
$ws.Range("J19").Value = 9.84506959824709
$ws.Range("K19").Value = 11.69111409903734
$ws.Range("L19").Value = 11.17332060155988
$ws.Range("N19").Value = 20.36932089331878
$ws.Range("O19").Value = 27.74006198349479

# Row 20
$ws.Range("B20").Value = 14.69405745391376
$ws.Range("C20").Value = 6.572885713387898
$ws.Range("D20").Value = 8.115781886088802
$ws.Range("E20").Value = 12.57509508140073
$ws.Range("F20").Value = 36.03528057671291
$ws.Range("H20").Value = 7.344005520526261
$ws.Range("I20").Value = 26.86605428694722
$ws.Range("J20").Value = 9.839048418214173
$ws.Range("K20").Value = 11.74975055807828
$ws.Range("L20").Value = 11.17641248555065
$ws.Range("N20").Value = 20.34889555615351
$ws.Range("O20").Value = 27.71702798370463

# Row 21
$ws.Range("B21").Value = 14.97142590725877
$ws.Range("C21").Value = 6.653852580864676
$ws.Range("D21").Value = 8.147538781344869
$ws.Range("E21").Value = 12.56487562939777
$ws.Range("F21").Value = 35.98119284280445
$ws.Range("H21").Value = 7.344005520526261
$ws.Range("I21").Value = 26.77872325483293
$ws.Range("J21").Value = 9.819595443961511
$ws.Range("K21").Value = 11.9463932347314
$ws.Range("L21").Value = 11.18836356357227
$ws.Range("N21").Value = 20.2821895711436
$ws.Range("O21").Value = 27.6451659290589

# Row 22
$ws.Range("B22").Value = 15.15165525765903
$ws.Range("C22").Value = 6.70577723186904
$ws.Range("D22").Value = 8.169129237766555
$ws.Range("E22").Value = 12.55947984903793
$ws.Range("F22").Value = 35.95179267687945
$ws.Range("H22").Value = 7.344005520526261
$ws.Range("I22").Value = 26.72525338307703
$ws.Range("J22").Value = 9.807459262639782
$ws.Range("K22").Value = 12.07444338213865
$ws.Range("L22").Value = 11.19732934977735
$ws.Range("N22").Value = 20.24001910613133
$ws.Range("O22").Value = 27.60232376598799

# Row 23
$ws.Range("B23").Value = 15.05557741817676
$ws.Range("C23").Value = 6.678157052358926
$ws.Range("D23").Value = 8.157530896347714
$ws.Range("E23").Value = 12.56224252841813
$ws.Range("F23").Value = 35.9669397528428
$ws.Range("H23").Value = 7.344005520526261
$ws.Range("I23").Value = 26.75346127858113
$ws.Range("J23").Value = 9.813883855887017
$ws.Range("K23").Value = 12.00615559000896
$ws.Range("L23").Value = 11.19243909205118
$ws.Range("N23").Value = 20.26239617267866
$ws.Range("O23").Value = 27.62481192294606

# Row 24
$ws.Range("B24").Value = 14.6895104491811
$ws.Range("C24").Value = 6.571546110272873
$ws.Range("D24").Value = 8.115277235202113
$ws.Range("E24").Value = 12.57528383428297
$ws.Range("F24").Value = 36.03626504304078
$ws.Range("H24").Value = 7.344005520526261
$ws.Range("I24").Value = 26.86754253733247
$ws.Range("J24").Value = 9.839376143263252
$ws.Range("K24").Value = 11.74653140917466
$ws.Range("L24").Value = 11.17623679089928
$ws.Range("N24").Value = 20.35000997954116
$ws.Range("O24").Value = 27.71827210858207

# Row 25
$ws.Range("B25").Value = 14.29204114661028
$ws.Range("C25").Value = 6.452520881846188
$ws.Range("D25").Value = 8.073417207532975
$ws.Range("E25").Value = 12.5948857126964
$ws.Range("F25").Value = 36.13657824024493
$ws.Range("H25").Value = 7.344005520526261
$ws.Range("I25").Value = 27.00607662572105
$ws.Range("J25").Value = 9.869365514692667
$ws.Range("K25").Value = 11.4657408058245
$ws.Range("L25").Value = 11.16378038829153
$ws.Range("N25").Value = 20.45067305953376
$ws.Range("O25").Value = 27.83676911625497
